# Formed the consolidated report
# Fill in the "Absent" (column H) values for rows where the attendance
# total is present but the "Absent" flag is still missing/blank, or needs
# to be corrected to reflect whether the student was actually marked
# present ("Real" = column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that need their Absent (H) value consolidated:
#   Row 3  -> Absent = 1 (Real=0, Invalid=1)
#   Row 7  -> Absent = 0 (Real=1)
#   Row 10 -> Absent = 1 (Real=0)
#   Row 11 -> Absent = 0 (Real=1)
#   Row 13 -> Absent = 1 (Real=0)
#   Row 14 -> Absent = 0 (Real=1)
#   Row 19 -> Absent = 1 (Real=0)
#   Row 20 -> Absent = 0 (Real=1)

$ws.Range("H3").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
